$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "blog" article (ser: 80) went live. The feed keeps its 3 most
# recent "blog" slots in row 7 (I7, E7, C7 - oldest to newest visually).
# The oldest one (ser: 77) drops off, the other two shift over, and the
# new article (ser: 80) takes the freshest slot.
$ws.Range("I7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 78"
$ws.Range("E7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 79"
$ws.Range("C7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 80"

# The author scrolled the sheet over to bring the new content into view
# and left the cursor on the (now) newly-shifted E7 cell.
$ws.Range("E7").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 4
